$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.703.36'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '3.322.75'
$ws.Range('E3').Value = '  +4.32%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.320.88'
$ws.Range('E8').Value = '  +4.37%  '
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('E10').Value = '  +2.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.54'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.67%  '
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.95'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').Value = '3.864.36'
$ws.Range('E15').Value = '  +4.53%  '
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '3.318.30'
$ws.Range('E17').Value = '  +4.74%  '
$ws.Range('D18').Value = '63.778.92'
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('E19').Value = '  +2.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '479.88'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.737'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.79'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('E27').Value = '  +2.54%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.20'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.17'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.98'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.107'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.10'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.36'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '0.0₃0745'
$ws.Range('E38').Value = '  +4.41%  '
$ws.Range('E39').Value = '  +3.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '435.37'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range('D41').Value = '3.100.68'
$ws.Range('E41').Value = '  +4.67%  '
$ws.Range('E42').Value = '  +7.07%  '
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.267'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.26'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '37.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +16.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.44'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.10%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.31'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.43%  '
$ws.Range('E51').Value = '  +0.16%  '
